$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "nazev"
$ws.Range("B1").Value = "zkratka"
$ws.Range("C1").Value = "prednasejiciUcitIdno"
$ws.Range("D1").Value = "jmena"

# Column A (nazev)
$ws.Cells.Item(2, 1).Value = 'Teoretická informatika a matematika'
$ws.Cells.Item(3, 1).Value = 'Paralelní programování'
$ws.Cells.Item(4, 1).Value = 'Paralelní programování'
$ws.Cells.Item(5, 1).Value = 'Matematika pro informatiky'
$ws.Cells.Item(6, 1).Value = 'Teorie formálních jazyků'
$ws.Cells.Item(7, 1).Value = 'Algoritmy a datové struktury'
$ws.Cells.Item(8, 1).Value = 'Algoritmy a datové struktury'
$ws.Cells.Item(9, 1).Value = 'System Simulation'
$ws.Cells.Item(10, 1).Value = 'System Simulation'
$ws.Cells.Item(11, 1).Value = 'Počítačové sítě a protokoly'
$ws.Cells.Item(12, 1).Value = 'Datová úložiště a nástroje pro Big Data'
$ws.Cells.Item(13, 1).Value = 'Datová úložiště a nástroje pro Big Data'
$ws.Cells.Item(14, 1).Value = 'Information Security'
$ws.Cells.Item(15, 1).Value = 'Automatické řízení'
$ws.Cells.Item(16, 1).Value = 'Teoretická informatika a matematika'
$ws.Cells.Item(17, 1).Value = 'Analogová elektronika'
$ws.Cells.Item(18, 1).Value = 'Analogová elektronika'
$ws.Cells.Item(19, 1).Value = 'Analogová elektronika'
$ws.Cells.Item(20, 1).Value = 'Analýza síťové komunikace'
$ws.Cells.Item(21, 1).Value = 'Python and R for Data Science'

# Column B (zkratka) -- "0177" must stay text, not become the number 177,
# so it goes through the text-formula + paste-as-values route too.
$ws.Cells.Item(2, 2).Value = 'YTIM'
$ws.Cells.Item(3, 2).Value = 'KPPG'
$ws.Cells.Item(4, 2).Value = 'PPG'
$ws.Cells.Item(5, 2).Value = 'MAI'
$ws.Cells.Item(6, 2).Value = 'TFL'
$ws.Cells.Item(7, 2).Value = 'DSA'
$ws.Cells.Item(8, 2).Value = 'KDSA'
$ws.Cells.Item(9, 2).Value = 'SYS'
$ws.Cells.Item(10, 2).Value = 'KSYS'
$ws.Cells.Item(11, 2).Value = 'PSP'
$ws.Cells.Item(12, 2).Value = 'BIG'
$ws.Cells.Item(13, 2).Value = 'KBIG'
$ws.Cells.Item(14, 2).Formula = "=""0177"""
$ws.Cells.Item(14, 2).Copy()
$ws.Cells.Item(14, 2).PasteSpecial(-4163)
$ws.Cells.Item(15, 2).Value = 'AUC'
$ws.Cells.Item(16, 2).Value = 'YTIM'
$ws.Cells.Item(17, 2).Value = 'ANE'
$ws.Cells.Item(18, 2).Value = 'KAEL'
$ws.Cells.Item(19, 2).Value = 'KANE'
$ws.Cells.Item(20, 2).Value = 'ASK'
$ws.Cells.Item(21, 2).Value = 'EPYR'
$excel.CutCopyMode = $false

# Column D (jmena) -- values carry a literal leading apostrophe in the source
# data, so they are written via a text formula + paste-as-values so the
# apostrophe is kept as literal text instead of Excel's quote-prefix marker.
$ws.Cells.Item(2, 4).Formula = "=""'doc. PaedDr. Petr Eisenmann CSc."""
$ws.Cells.Item(2, 4).Copy()
$ws.Cells.Item(2, 4).PasteSpecial(-4163)
$ws.Cells.Item(2, 4).VerticalAlignment = -4108
$ws.Cells.Item(3, 4).Formula = "=""'prof. Ing. Martin Lísal DSc."""
$ws.Cells.Item(3, 4).Copy()
$ws.Cells.Item(3, 4).PasteSpecial(-4163)
$ws.Cells.Item(3, 4).VerticalAlignment = -4108
$ws.Cells.Item(4, 4).Formula = "=""'prof. Ing. Martin Lísal DSc."""
$ws.Cells.Item(4, 4).Copy()
$ws.Cells.Item(4, 4).PasteSpecial(-4163)
$ws.Cells.Item(4, 4).VerticalAlignment = -4108
$ws.Cells.Item(5, 4).Formula = "=""'doc. RNDr. Jiří Felcman CSc."""
$ws.Cells.Item(5, 4).Copy()
$ws.Cells.Item(5, 4).PasteSpecial(-4163)
$ws.Cells.Item(5, 4).VerticalAlignment = -4108
$ws.Cells.Item(6, 4).Formula = "=""'doc. RNDr. Karel Oliva Dr.'"""
$ws.Cells.Item(6, 4).Copy()
$ws.Cells.Item(6, 4).PasteSpecial(-4163)
$ws.Cells.Item(6, 4).VerticalAlignment = -4108
$ws.Cells.Item(7, 4).Formula = "=""'doc. RNDr. Karel Oliva Dr."""
$ws.Cells.Item(7, 4).Copy()
$ws.Cells.Item(7, 4).PasteSpecial(-4163)
$ws.Cells.Item(7, 4).VerticalAlignment = -4108
$ws.Cells.Item(8, 4).Formula = "=""'doc. RNDr. Karel Oliva Dr."""
$ws.Cells.Item(8, 4).Copy()
$ws.Cells.Item(8, 4).PasteSpecial(-4163)
$ws.Cells.Item(8, 4).VerticalAlignment = -4108
$ws.Cells.Item(9, 4).Formula = "=""'doc. RNDr. Jan Jirsák Ph.D."""
$ws.Cells.Item(9, 4).Copy()
$ws.Cells.Item(9, 4).PasteSpecial(-4163)
$ws.Cells.Item(9, 4).VerticalAlignment = -4108
$ws.Cells.Item(10, 4).Formula = "=""'doc. RNDr. Jan Jirsák Ph.D."""
$ws.Cells.Item(10, 4).Copy()
$ws.Cells.Item(10, 4).PasteSpecial(-4163)
$ws.Cells.Item(10, 4).VerticalAlignment = -4108
$ws.Cells.Item(11, 4).Formula = "=""'Ing. Toni Koluch Ph.D.'"""
$ws.Cells.Item(11, 4).Copy()
$ws.Cells.Item(11, 4).PasteSpecial(-4163)
$ws.Cells.Item(11, 4).VerticalAlignment = -4108
$ws.Cells.Item(12, 4).Formula = "=""'Ing. Václav Valenta'"""
$ws.Cells.Item(12, 4).Copy()
$ws.Cells.Item(12, 4).PasteSpecial(-4163)
$ws.Cells.Item(12, 4).VerticalAlignment = -4108
$ws.Cells.Item(13, 4).Formula = "=""'Ing. Václav Valenta'"""
$ws.Cells.Item(13, 4).Copy()
$ws.Cells.Item(13, 4).PasteSpecial(-4163)
$ws.Cells.Item(13, 4).VerticalAlignment = -4108
$ws.Cells.Item(14, 4).Formula = "=""'Hoon Ko Ph.D.'"""
$ws.Cells.Item(14, 4).Copy()
$ws.Cells.Item(14, 4).PasteSpecial(-4163)
$ws.Cells.Item(14, 4).VerticalAlignment = -4108
$ws.Cells.Item(15, 4).Formula = "=""'doc. Ing. Mgr. Petr Klán CSc.'"""
$ws.Cells.Item(15, 4).Copy()
$ws.Cells.Item(15, 4).PasteSpecial(-4163)
$ws.Cells.Item(15, 4).VerticalAlignment = -4108
$ws.Cells.Item(16, 4).Formula = "=""'RNDr. Veronika Pitrová PhD., Ph.D."""
$ws.Cells.Item(16, 4).Copy()
$ws.Cells.Item(16, 4).PasteSpecial(-4163)
$ws.Cells.Item(16, 4).VerticalAlignment = -4108
$ws.Cells.Item(17, 4).Formula = "=""'doc. RNDr. František Lustig CSc.'"""
$ws.Cells.Item(17, 4).Copy()
$ws.Cells.Item(17, 4).PasteSpecial(-4163)
$ws.Cells.Item(17, 4).VerticalAlignment = -4108
$ws.Cells.Item(18, 4).Formula = "=""'doc. RNDr. František Lustig CSc.'"""
$ws.Cells.Item(18, 4).Copy()
$ws.Cells.Item(18, 4).PasteSpecial(-4163)
$ws.Cells.Item(18, 4).VerticalAlignment = -4108
$ws.Cells.Item(19, 4).Formula = "=""'doc. RNDr. František Lustig CSc.'"""
$ws.Cells.Item(19, 4).Copy()
$ws.Cells.Item(19, 4).PasteSpecial(-4163)
$ws.Cells.Item(19, 4).VerticalAlignment = -4108
$ws.Cells.Item(20, 4).Formula = "=""'Ing. Vojtěch Šindler'"""
$ws.Cells.Item(20, 4).Copy()
$ws.Cells.Item(20, 4).PasteSpecial(-4163)
$ws.Cells.Item(20, 4).VerticalAlignment = -4108
$ws.Cells.Item(21, 4).Formula = "=""'Ricardo Rodriguez Jorge Ph.D."""
$ws.Cells.Item(21, 4).Copy()
$ws.Cells.Item(21, 4).PasteSpecial(-4163)
$ws.Cells.Item(21, 4).VerticalAlignment = -4108
$excel.CutCopyMode = $false

# Column C (prednasejiciUcitIdno, numeric)
$ws.Cells.Item(2, 3).Value = 261
$ws.Cells.Item(3, 3).Value = 609
$ws.Cells.Item(4, 3).Value = 609
$ws.Cells.Item(5, 3).Value = 816
$ws.Cells.Item(6, 3).Value = 2240
$ws.Cells.Item(7, 3).Value = 2240
$ws.Cells.Item(8, 3).Value = 2240
$ws.Cells.Item(9, 3).Value = 2855
$ws.Cells.Item(10, 3).Value = 2855
$ws.Cells.Item(11, 3).Value = 3458
$ws.Cells.Item(12, 3).Value = 4195
$ws.Cells.Item(13, 3).Value = 4195
$ws.Cells.Item(14, 3).Value = 4685
$ws.Cells.Item(15, 3).Value = 5039
$ws.Cells.Item(16, 3).Value = 5205
$ws.Cells.Item(17, 3).Value = 5232
$ws.Cells.Item(18, 3).Value = 5232
$ws.Cells.Item(19, 3).Value = 5232
$ws.Cells.Item(20, 3).Value = 7152
$ws.Cells.Item(21, 3).Value = 8514
